$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A="4A";  B=0.853},
    @{Row=3;  A="4P";  B=0.893},
    @{Row=4;  A="4S";  B=0.9508333329999999},
    @{Row=5;  A="5P";  B=0.99025},
    @{Row=6;  A="5W";  B=0.873},
    @{Row=7;  A="8F";  B=0.8784999999999999},
    @{Row=8;  A="8H";  B=0.9105},
    @{Row=9;  A="8K";  B=0.9385},
    @{Row=10; A="10H"; B=0.9025},
    @{Row=11; A="20H"; B=0.9695},
    @{Row=12; A="55I"; B=0.916},
    @{Row=13; A="71K"; B=0.9895}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 1).Value = $item.A
    $ws.Cells.Item($item.Row, 2).Value = $item.B
}
